# Daily attendance processing - 2025-11-04 18:29:47
# Normalize the "Recorded By" (column G) entries: move the signed-in
# reviewer's e-mail ("dnasr281@gmail.com") to the front of the list,
# and for the rows that double-logged the "System" account (once as
# "System" and once as "system"), swap those two tokens so the
# lowercase "system" entry comes first. Rows that don't contain
# "dnasr281@gmail.com" and don't have the duplicated System/system
# pair are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$colG = 7
$colStatus = 9
$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Column I ("Status") is always populated ("Recorded" vs "Pending").
    # Check it first so we never have to read/touch column G on rows
    # where it is genuinely blank (merely reading a blank cell in this
    # engine causes it to be serialized back out as an empty string).
    $status = $ws.Cells.Item($r, $colStatus).Value2
    if ($null -eq $status -or -not $status.Equals("Recorded")) {
        continue
    }

    $cell = $ws.Cells.Item($r, $colG)
    $value = $cell.Value2

    if ($null -eq $value) {
        continue
    }
    if ($value -notlike "*,*") {
        continue
    }

    $parts = $value -split ", "

    # NOTE: the "-eq"/"-ceq" string operators in this engine are
    # case-insensitive, so use the .Equals() instance method (which is
    # ordinal / case-sensitive) to tell "System" and "system" apart.
    $dnasrIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i].Equals("dnasr281@gmail.com")) {
            $dnasrIndex = $i
        }
    }

    $newParts = $null

    if ($dnasrIndex -ge 0) {
        # Move the dnasr281@gmail.com entry to the front, preserving the
        # relative order of the remaining entries.
        $newParts = @($parts[$dnasrIndex])
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $dnasrIndex) {
                $newParts += $parts[$i]
            }
        }
    }
    elseif ($parts.Length -ge 2 -and $parts[0].Equals("System") -and $parts[1].Equals("system")) {
        # Duplicate System/system entry - swap the first two tokens,
        # keep everything after them (e.g. backup@backdoor.com) in place.
        $newParts = @($parts[1], $parts[0])
        for ($i = 2; $i -lt $parts.Length; $i++) {
            $newParts += $parts[$i]
        }
    }

    if ($null -ne $newParts) {
        $newValue = $newParts -join ", "
        # Use .Equals() (ordinal/case-sensitive) instead of -eq/-ne,
        # since -eq/-ne are case-insensitive in this engine and would
        # treat "System, system, X" and "system, System, X" as equal.
        if (-not $newValue.Equals($value)) {
            $cell.Value = $newValue
            $changed++
        }
    }
}

Write-Output "Updated $changed cell(s) in column G"
